$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
# Columns D..L are new/renamed headers; M..O shift from former E..G
$ws.Range("A1").Value = "Best Estimator"
$ws.Range("B1").Value = "Best Score"
$ws.Range("C1").Value = "Best Params"
$ws.Range("D1").Value = "CV Train F1"
$ws.Range("E1").Value = "CV Test F1"
$ws.Range("F1").Value = "Validation F1"
$ws.Range("G1").Value = "CV Train Precision"
$ws.Range("H1").Value = "CV Test Precision"
$ws.Range("I1").Value = "Validation Precision"
$ws.Range("J1").Value = "CV Train Recall"
$ws.Range("K1").Value = "CV Test Recall"
$ws.Range("L1").Value = "Validation Recall"
$ws.Range("M1").Value = "Y Val (Validation)"
$ws.Range("N1").Value = "Y Pred (Validation)"
$ws.Range("O1").Value = "Seed"

# Ensure the new header cells (D1:O1) use the same bold/bordered style as the original header
$ws.Range("A1").Copy()
$ws.Range("D1:O1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 2 ---
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', MinMaxScaler()), ('selector', None),`n                ('model',`n                 LogisticRegression(class_weight='balanced', l1_ratio=0.95,`n                                    max_iter=1000, penalty='elasticnet',`n                                    random_state=42, solver='saga'))])"
$ws.Range("B2").Value = 0.6004761904761905
$ws.Range("C2").Value = "{'selector': None, 'scaler': MinMaxScaler(), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.95, 'model__class_weight': 'balanced'}"
$ws.Range("D2").Value = 0.7060911715132402
$ws.Range("E2").Value = 0.5508879546379546
$ws.Range("F2").Value = 0.7727272727272727
$ws.Range("G2").Value = 0.658559013431942
$ws.Range("H2").Value = 0.5475314153439154
$ws.Range("I2").Value = 0.7727272727272727
$ws.Range("J2").Value = 0.7962301587301587
$ws.Range("K2").Value = 0.6229166666666667
$ws.Range("L2").Value = 0.7727272727272727
$ws.Range("M2").Value = "[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1 0 1 1 1 1 0 0 0 0 1 0 1]"
$ws.Range("N2").Value = "[1 0 0 1 1 1 0 0 1 0 0 1 1 1 0 0 0 1 1 1 1 1 1 1 0 1 1 0 1 0 1 0 0 1 1 1]"
$ws.Range("O2").Value = 42

# --- Row 3 ---
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fd58723c160>),`n                ('model',`n                 LogisticRegression(class_weight='balanced', l1_ratio=0.01,`n                                    max_iter=1000, penalty='elasticnet',`n                                    random_state=42, solver='saga'))])"
$ws.Range("B3").Value = 0.5771428571428572
$ws.Range("C3").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd58722af10>, 'scaler': MinMaxScaler(), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.01, 'model__class_weight': 'balanced'}"
$ws.Range("D3").Value = 0.6965802274313946
$ws.Range("E3").Value = 0.5314406658156657
$ws.Range("F3").Value = 0.6956521739130435
$ws.Range("G3").Value = 0.68167516409576
$ws.Range("H3").Value = 0.4948495370370371
$ws.Range("I3").Value = 0.7272727272727273
$ws.Range("J3").Value = 0.7388020833333333
$ws.Range("K3").Value = 0.6083333333333333
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 0 1 0 1 1 0 1 1 0 1 1 1]"
$ws.Range("N3").Value = "[0 1 1 1 1 0 1 1 1 1 0 0 1 0 1 1 1 1 1 0 1 0 1 1 1 1 0 0 0 0 0 0 0 1 1 1]"
$ws.Range("O3").Value = 69

# --- Row 4 ---
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fd58722af40>),`n                ('model',`n                 LogisticRegression(class_weight='balanced', l1_ratio=0.1,`n                                    max_iter=1000, penalty='elasticnet',`n                                    random_state=42, solver='saga'))])"
$ws.Range("B4").Value = 0.5695238095238095
$ws.Range("C4").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd586fb0fd0>, 'scaler': MinMaxScaler(), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.1, 'model__class_weight': 'balanced'}"
$ws.Range("D4").Value = 0.6861869258981556
$ws.Range("E4").Value = 0.4950939454064454
$ws.Range("F4").Value = 0.5909090909090908
$ws.Range("G4").Value = 0.6841178586910491
$ws.Range("H4").Value = 0.5054414682539683
$ws.Range("I4").Value = 0.7222222222222222
$ws.Range("J4").Value = 0.6947916666666667
$ws.Range("K4").Value = 0.5091666666666667
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = "[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1 0 1 0 1 0 1 0 1 1 1 0 1]"
$ws.Range("N4").Value = "[0 1 1 0 1 1 0 1 1 0 0 0 0 0 1 0 0 1 0 1 0 0 1 0 1 1 0 1 1 1 1 0 1 0 1 0]"
$ws.Range("O4").Value = 23

# --- Row 5 ---
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fd586fb0a00>),`n                ('model',`n                 LogisticRegression(class_weight='balanced', l1_ratio=0.95,`n                                    max_iter=1000, penalty='elasticnet',`n                                    random_state=42, solver='saga'))])"
$ws.Range("B5").Value = 0.7808333333333333
$ws.Range("C5").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd5871d72e0>, 'scaler': MinMaxScaler(), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.95, 'model__class_weight': 'balanced'}"
$ws.Range("D5").Value = 0.7079136418993711
$ws.Range("E5").Value = 0.6032895114145114
$ws.Range("F5").Value = 0.55
$ws.Range("G5").Value = 0.6682680302774188
$ws.Range("H5").Value = 0.5695122354497354
$ws.Range("I5").Value = 0.6111111111111112
$ws.Range("J5").Value = 0.7925595238095239
$ws.Range("K5").Value = 0.7004166666666666
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = "[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0 1 1 1 1 1 0]"
$ws.Range("N5").Value = "[0 1 1 0 1 0 0 1 0 0 1 1 0 0 1 0 0 0 0 0 0 1 1 1 1 0 1 1 0 1 0 1 1 1 0 1]"
$ws.Range("O5").Value = 99

# --- Row 6 ---
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector',`n                 <__main__.NamedFeatureSelector object at 0x7fd5871d7280>),`n                ('model',`n                 LogisticRegression(class_weight='balanced', l1_ratio=0.5,`n                                    max_iter=1000, penalty='elasticnet',`n                                    random_state=42, solver='saga'))])"
$ws.Range("B6").Value = 0.6866666666666668
$ws.Range("C6").Value = "{'selector': <__main__.NamedFeatureSelector object at 0x7fd5871d7670>, 'scaler': MinMaxScaler(), 'model__solver': 'saga', 'model__penalty': 'elasticnet', 'model__l1_ratio': 0.5, 'model__class_weight': 'balanced'}"
$ws.Range("D6").Value = 0.775558231363358
$ws.Range("E6").Value = 0.6353811466311465
$ws.Range("F6").Value = 0.5714285714285713
$ws.Range("G6").Value = 0.740178670251466
$ws.Range("H6").Value = 0.5693716931216931
$ws.Range("I6").Value = 0.5454545454545454
$ws.Range("J6").Value = 0.8401515151515152
$ws.Range("K6").Value = 0.75
$ws.Range("L6").Value = 0.6
$ws.Range("M6").Value = "[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1 1 0 1 0 1 1 1 1 1 1 1 0]"
$ws.Range("N6").Value = "[1 1 1 0 1 1 1 0 0 0 0 1 0 1 1 0 1 1 1 1 1 1 1 1 0 1 1 0 0 1 1 0 1 0 0 0]"
$ws.Range("O6").Value = 89

# Re-fit row heights so the newly entered multi-line "Best Estimator" text
# doesn't leave a stray custom row height behind (matches original formatting)
$ws.Rows.Item(2).AutoFit()
$ws.Rows.Item(3).AutoFit()
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()
$ws.Rows.Item(6).AutoFit()
